$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Product$ProductBuilder block (rows 2-5): reorder fields so
# availableItems/reservedItems come before id/name
$ws.Range("B2").Value = "availableItems"
$ws.Range("D2").Value = "int"

$ws.Range("B3").Value = "reservedItems"
$ws.Range("D3").Value = "int"

$ws.Range("B4").Value = "id"
$ws.Range("D4").Value = "java.lang.Long"

$ws.Range("B5").Value = "name"
$ws.Range("D5").Value = "java.lang.String"

# Product block (rows 8-11): same reordering
$ws.Range("B8").Value = "availableItems"
$ws.Range("D8").Value = "int"

$ws.Range("B9").Value = "reservedItems"
$ws.Range("D9").Value = "int"

$ws.Range("B10").Value = "id"
$ws.Range("D10").Value = "java.lang.Long"

$ws.Range("B11").Value = "name"
$ws.Range("D11").Value = "java.lang.String"

# OrderServiceImpl block (rows 14-16): move template field earlier,
# push SOURCE field down
$ws.Range("B14").Value = "template"
$ws.Range("D14").Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Range("B15").Value = "log"
$ws.Range("D15").Value = "org.slf4j.Logger"

$ws.Range("B16").Value = "SOURCE"
$ws.Range("D16").Value = "domain.OrderSource"
